# Updated symbol list with refreshed Price (D) and Volume(1h) (E) values.
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the workbook's original text-formatted Price/Volume columns)
# instead of auto-converting numeric-looking strings into real numbers;
# the Style reset afterwards clears the quote-prefix formatting that the
# apostrophe trick would otherwise leave behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.70%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'29.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.36%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.190"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.05736"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.01%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.570"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.32%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8590"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.8642"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.01%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1365"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.56%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.07077"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.75%"
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'6.40%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.09377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.21%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.001538"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.91%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.0006012"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-94.10%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.006031"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.90%"
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'5,228.21%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.495"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.103"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.11%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.186"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.3204"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.74%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.03307"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.49%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.1290"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.30%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'3.477"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.18%"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'1.70%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.1380"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'0.001227"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'1.01%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.004992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'11.50%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.0001210"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'2.59%"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.03752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1071"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.30%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'0.87%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.003521"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-41.01%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.008462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-11.74%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005288"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.24%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.10%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.05702"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-43.51%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.002258"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-10.75%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
